# Adds header labels (A1) to the data tables that feed the charts, fixes
# accentuation of Portuguese labels, updates a couple of values, and removes
# the now-unused "Teto" row on the emissions sheet.
#
# xlPasteFormats (used below) = -4122. Copy/PasteSpecial is used (instead of
# a plain `.Style =` assignment, which only touches the named cell style and
# not the direct border/font/alignment formatting) so that new header cells
# pick up exactly the same direct formatting as their neighbours, and label
# cells that lose their bold/bordered look fall back to the plain default
# format -- matching what Excel itself would do when you copy formats around.
$xlPasteFormats = -4122

function Set-CellFormatLike {
    param($TargetRange, $SourceRange)
    $SourceRange.Copy() | Out-Null
    $TargetRange.PasteSpecial($xlPasteFormats) | Out-Null
}

$wb = $excel.ActiveWorkbook

# --- Sheets 1-4: "Fonte/Tecnologia" header + accent fixes ---------------
$fonteSheets = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($sheetName in $fonteSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # New header for column A, styled like the other header cells (B1:E1)
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    Set-CellFormatLike $ws.Range("A1") $ws.Range("B1")

    # Fix accentuation / abbreviations
    $ws.Range("A2").Value = "Hidro"
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A5").Value = "Nuclear"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A7").Value = "Biomassa"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A9").Value = "Solar"
    $ws.Range("A10").Value = "Outros"
    $ws.Range("A11").Value = "Pot. Compl."
    $ws.Range("A12").Value = "GD"

    # These label cells drop the bold/bordered header style -> plain format
    for ($r = 2; $r -le 12; $r++) {
        Set-CellFormatLike $ws.Cells.Item($r, 1) $ws.Cells.Item($r, 2)
    }
}

# --- Sheet 5: "Emissoes Totais (MtCO2eq)" --------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("A1").Value = "Período"
Set-CellFormatLike $ws5.Range("A1") $ws5.Range("B1")

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
Set-CellFormatLike $ws5.Cells.Item(2, 1) $ws5.Cells.Item(2, 2)
Set-CellFormatLike $ws5.Cells.Item(3, 1) $ws5.Cells.Item(3, 2)

# Remove the "Teto" row entirely
$ws5.Rows.Item(4).Delete()

# --- Sheet 6: "Custo Total (bilhões de R$)" ------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws6.Range("A1").Value = "Tipo Expansão"
Set-CellFormatLike $ws6.Range("A1") $ws6.Range("B1")

# B1 keeps its header style but its text changes from "Custo" to "2015" --
# force it to stay text (like every other sheet's B1) instead of letting it
# be auto-recognised as a number. Setting NumberFormat first borrows a fresh
# style slot, so re-stamp the header format (copied from A1, set up above)
# afterwards to land back on the shared header style.
$ws6.Range("B1").NumberFormat = "@"
$ws6.Range("B1").Value = "2015"
Set-CellFormatLike $ws6.Range("B1") $ws6.Range("A1")

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A3").Value = "Expansão por GD"
Set-CellFormatLike $ws6.Cells.Item(2, 1) $ws6.Cells.Item(2, 2)
Set-CellFormatLike $ws6.Cells.Item(3, 1) $ws6.Cells.Item(3, 2)

$ws6.Range("B2").Value = 569
$ws6.Range("B3").Value = 99
